# Edit script: apply the changes described by the commit
# 1. Fix several "name_en" hero-name strings on the name_heroes sheet
# 2. Insert a new worksheet "Sheet7" before "Sheet6" containing the full
#    alphabetical list of hero names (used for tooltip lookups / data validation)
# 3. Leave the existing "Sheet6" sheet content untouched
# 4. Re-select the first sheet (name_heroes) as the active sheet/tab

$wb = $excel.ActiveWorkbook

$replacements = @(
  @("`"name_en`": `"Corrin(M)`",", "`"name_en`": `"Corrin (M)`","),
  @("`"name_en`": `"Lonqu`",", "`"name_en`": `"Lon'qu`","),
  @("`"name_en`": `"Tiki(A)`",", "`"name_en`": `"Tiki (Adult)`","),
  @("`"name_en`": `"Tiki(Y)`",", "`"name_en`": `"Tiki (Young)`","),
  @("`"name_en`": `"Robin(F)`",", "`"name_en`": `"Robin (F)`","),
  @("`"name_en`": `"Corrin(F)`",", "`"name_en`": `"Corrin (F)`","),
  @("`"name_en`": `"Robin(M)`",", "`"name_en`": `"Robin (M)`","),
  @("`"name_en`": `"Lucina (Spring)`",", "`"name_en`": `"Lucina (Spring Festival)`","),
  @("`"name_en`": `"Camilla (Spring)`",", "`"name_en`": `"Camilla (Spring Festival)`","),
  @("`"name_en`": `"Xander (Spring)`",", "`"name_en`": `"Xander (Spring Festival)`","),
  @("`"name_en`": `"Chrom (Spring)`",", "`"name_en`": `"Chrom (Spring Festival)`","),
  @("`"name_en`": `"Caeda (Bride)`",", "`"name_en`": `"Caeda (Bridal Blessings)`","),
  @("`"name_en`": `"Charlotte (Bride)`",", "`"name_en`": `"Charlotte (Bridal Blessings)`","),
  @("`"name_en`": `"Cordelia (Bride)`",", "`"name_en`": `"Cordelia (Bridal Blessings)`","),
  @("`"name_en`": `"Lyn (Bride)`",", "`"name_en`": `"Lyn (Bridal Blessings)`","),
  @("`"name_en`": `"Frederick (Summer)`",", "`"name_en`": `"Frederick (Ylissean Summer)`","),
  @("`"name_en`": `"Gaius (Summer)`",", "`"name_en`": `"Gaius (Ylissean Summer)`","),
  @("`"name_en`": `"Robin(F) (Summer)`",", "`"name_en`": `"Robin (F) (Ylissean Summer)`","),
  @("`"name_en`": `"Tiki(A) (Summer)`",", "`"name_en`": `"Tiki (Adult) (Ylissean Summer)`","),
  @("`"name_en`": `"Corrin(F) (Summer)`",", "`"name_en`": `"Corrin (F) (Nohrian Summer)`","),
  @("`"name_en`": `"Leo (Summer)`",", "`"name_en`": `"Leo (Nohrian Summer)`","),
  @("`"name_en`": `"Xander (Summer)`",", "`"name_en`": `"Xander (Nohrian Summer)`","),
  @("`"name_en`": `"Elise (Summer)`",", "`"name_en`": `"Elise (Nohrian Summer)`",")
)
$nameHeroes = $wb.Worksheets.Item("name_heroes")
foreach ($pair in $replacements) {
  $nameHeroes.Cells.Replace($pair[0], $pair[1])
}

$sheet7Values = @(
  "Abel",
  "Alfonse",
  "Alm",
  "Amelia",
  "Anna",
  "Arden",
  "Arthur",
  "Arvis",
  "Athena",
  "Ayra",
  "Azama",
  "Azura (Performing Arts)",
  "Azura",
  "Barst",
  "Bartre",
  "Berkut",
  "Beruka",
  "Black Knight",
  "Boey",
  "Bruno",
  "Caeda (Bridal Blessings)",
  "Caeda",
  "Cain",
  "Camilla",
  "Camilla (Spring Festival)",
  "Camus",
  "Catria",
  "Cecilia",
  "Celica",
  "Charlotte (Bridal Blessings)",
  "Cherche",
  "Chrom",
  "Chrom (Spring Festival)",
  "Clair",
  "Clarine",
  "Clarisse",
  "Clive",
  "Cordelia (Bridal Blessings)",
  "Cordelia",
  "Corrin (F) (Nohrian Summer)",
  "Corrin (F)",
  "Corrin (M)",
  "Deirdre",
  "Delthea",
  "Donnel",
  "Dorcas",
  "Draug",
  "Effie",
  "Eirika",
  "Eldigan",
  "Elincia",
  "Elise (Nohrian Summer)",
  "Elise",
  "Eliwood",
  "Ephraim",
  "Est",
  "Fae",
  "Faye",
  "Felicia",
  "Fir",
  "Fjorm",
  "Florina",
  "Frederick",
  "Frederick (Ylissean Summer)",
  "Gaius",
  "Gaius (Ylissean Summer)",
  "generic_armored_axe",
  "generic_armored_lance",
  "generic_armored_sword",
  "generic_cavalry_axe",
  "generic_cavalry_bluetome",
  "generic_cavalry_bow",
  "generic_cavalry_greentome",
  "generic_cavalry_lance",
  "generic_cavalry_redtome",
  "generic_cavalry_staff",
  "generic_cavalry_sword",
  "generic_flying_axe",
  "generic_flying_lance",
  "generic_flying_sword",
  "generic_infantry_axe",
  "generic_infantry_bluetome",
  "generic_infantry_bow",
  "generic_infantry_dagger",
  "generic_infantry_greentome",
  "generic_infantry_lance",
  "generic_infantry_redtome",
  "generic_infantry_staff",
  "generic_infantry_sword",
  "Genny",
  "Gordin",
  "Gray",
  "Gunter",
  "Gwendolyn",
  "Hana",
  "Hawkeye",
  "Hector",
  "Henry",
  "Henry (Trick or Defeat!)",
  "Henry (Trick or Defeat)",
  "Hinata",
  "Hinoka",
  "Ike (Brave Heroes)",
  "Ike",
  "Inigo (Performing Arts)",
  "Innes",
  "Jaffar",
  "Jagen",
  "Jakob",
  "Jakob (Trick or Defeat!)",
  "Jakob (Trick or Defeat)",
  "Jeorge",
  "Joshua",
  "Julia",
  "Kagero",
  "Karel",
  "Katarina",
  "Klein",
  "Lachesis",
  "Lævateinn",
  "Laslow",
  "Legion",
  "Leo (Nohrian Summer)",
  "Leon",
  "Leo",
  "Lilina",
  "Linde",
  "Lissa",
  "Lloyd",
  "Loki",
  "Lon'qu",
  "Lucina (Brave Heroes)",
  "Lucina",
  "Lucina (Spring Festival)",
  "Lucius",
  "Lukas",
  "Luke",
  "Lute",
  "Lyn (Brave Heroes)",
  "Lyn (Bridal Blessings)",
  "Lyn",
  "Mae",
  "Maria",
  "Marth (Masked)",
  "Marth",
  "Mathilda",
  "Matthew",
  "Merric",
  "Mia",
  "Michalis",
  "Minerva",
  "Mist",
  "Narcian",
  "Navarre",
  "Nephenee",
  "Niles",
  "Ninian",
  "Nino",
  "nohero",
  "Nowi",
  "Nowi (Trick or Defeat!)",
  "Nowi (Trick or Defeat)",
  "Oboro",
  "Odin",
  "Ogma",
  "Olivia (Performing Arts)",
  "Olivia",
  "Olwen",
  "Oscar",
  "Palla",
  "Peri",
  "Priscilla",
  "Raigh",
  "Raven",
  "Rebecca",
  "Reinhardt",
  "Robin (F)",
  "Robin (F) (Ylissean Summer)",
  "Robin (M)",
  "Roderick",
  "Roy (Brave Heroes)",
  "Roy",
  "Ryoma",
  "Saber",
  "Saizo",
  "Sakura",
  "Sakura (Trick or Defeat!)",
  "Sakura (Trick or Defeat)",
  "Sanaki",
  "Selena",
  "Seliph",
  "Serra",
  "Seth",
  "Setsuna",
  "Shanna",
  "Sharena",
  "Sheena",
  "Shigure (Performing Arts)",
  "Sigurd",
  "Sonya",
  "Sophia",
  "Soren",
  "Stahl",
  "Subaki",
  "Sully",
  "Surtr",
  "Tailtiu",
  "Takumi",
  "Tana",
  "Tharja",
  "Tiki (Adult)",
  "Tiki (Adult) (Ylissean Summer)",
  "Tiki (Young)",
  "Titania",
  "Tobin",
  "Ursula",
  "Valter",
  "Veronica",
  "Virion",
  "Wrys",
  "Xander (Nohrian Summer)",
  "Xander",
  "Xander (Spring Festival)",
  "Zephiel"
)
$sheet6 = $wb.Worksheets.Item("Sheet6")
$sheet7 = $wb.Worksheets.Add($sheet6)
$sheet7.Name = "Sheet7"
for ($i = 0; $i -lt $sheet7Values.Count; $i++) {
  $sheet7.Cells.Item($i + 1, 1).Value = $sheet7Values[$i]
}
$sheet7.Columns.Item(1).ColumnWidth = 26.81

$sheet7.Activate()
$sheet7.Range("D213").Select()
$excel.ActiveWindow.ScrollRow = 202

$nameHeroes.Activate()
$nameHeroes.Range("B53").Select()
$excel.ActiveWindow.ScrollRow = 163
